$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns right before the "general_college_subjects.arts" column (R),
# shifting it (and everything after it) three columns to the right.
$ws.Range("R1:T1").EntireColumn.Insert()

# New header labels for the freshly inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data values (row 2) for those same columns.
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Normalize casing of a handful of row-2 text values.
$ws.Range("D2").Value = "considered"
$ws.Range("E2").Value = "considered"
$ws.Range("F2").Value = "considered"
$ws.Range("G2").Value = "considered"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "considered"
$ws.Range("J2").Value = "considered"

Write-Output "done"
